$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry (avoid Excel auto-converting numeric-looking strings)
# by temporarily setting NumberFormat to Text, then clearing formats
# afterwards so styling matches the original (unstyled) cells exactly.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.150.04'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.566.96'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.46'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.00'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.75%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.60'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.024.65'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.065.42'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000146'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.574.12'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '340.91'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.84'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.40%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.93'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.681.99'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.64'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.91%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.48'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.48'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.89'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +8.83%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0824'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '177.27'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '436.01'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.44%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.405'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.28'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.50'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.62%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '152.01'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.11'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0552'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.607'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.38'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.73'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.30%  '

# Strip the temporary Text number-format so styling reverts to default
$ws.Range('B2:E51').ClearFormats()

